# Auto-generated PowerShell Excel COM-interop script
# Applies the numeric cell updates described by the commit diff
# for '上海-漫展信息.xlsx' (Shanghai comic/anime convention info)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1620
$ws.Range("F6").Value = 796
$ws.Range("G6").Value = 39.9
$ws.Range("F7").Value = 699
$ws.Range("F8").Value = 1288
$ws.Range("F9").Value = 2617
$ws.Range("F10").Value = 1346
$ws.Range("F11").Value = 574
$ws.Range("F12").Value = 2317
$ws.Range("F13").Value = 2042
$ws.Range("F14").Value = 717
$ws.Range("F15").Value = 6481
$ws.Range("F16").Value = 123
$ws.Range("F17").Value = 1228
$ws.Range("F18").Value = 140
$ws.Range("F19").Value = 1477
$ws.Range("F20").Value = 1331
$ws.Range("F21").Value = 1194
$ws.Range("F23").Value = 2296
$ws.Range("F25").Value = 731
$ws.Range("F26").Value = 242
$ws.Range("F27").Value = 5303
$ws.Range("F28").Value = 286
$ws.Range("F29").Value = 1254
$ws.Range("F30").Value = 46
$ws.Range("F31").Value = 3711
$ws.Range("F33").Value = 1687
$ws.Range("F35").Value = 159
$ws.Range("F36").Value = 274
$ws.Range("F39").Value = 392
$ws.Range("F40").Value = 1768
$ws.Range("F42").Value = 105
$ws.Range("F43").Value = 899
$ws.Range("F45").Value = 513
$ws.Range("F46").Value = 47
$ws.Range("F48").Value = 60
$ws.Range("F49").Value = 78

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 441
$ws.Range("F10").Value = 15
$ws.Range("F11").Value = 389
$ws.Range("F13").Value = 138
$ws.Range("F15").Value = 965
$ws.Range("F22").Value = 242
$ws.Range("F23").Value = 356
$ws.Range("F26").Value = 81
$ws.Range("F27").Value = 81
$ws.Range("F30").Value = 301
$ws.Range("F31").Value = 39
$ws.Range("F35").Value = 46
$ws.Range("G35").Value = 149
$ws.Range("F36").Value = 111
$ws.Range("F38").Value = 197

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3297
$ws.Range("F5").Value = 397
$ws.Range("F8").Value = 765
$ws.Range("F10").Value = 2786
$ws.Range("F11").Value = 288
$ws.Range("F12").Value = 530
$ws.Range("F13").Value = 570
$ws.Range("F14").Value = 1160

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 765
$ws.Range("F6").Value = 2786
$ws.Range("F7").Value = 1620
$ws.Range("F8").Value = 796
$ws.Range("G8").Value = 39.9
$ws.Range("F9").Value = 699
$ws.Range("F10").Value = 1288
$ws.Range("F11").Value = 2617
$ws.Range("F12").Value = 1346
$ws.Range("F13").Value = 574
$ws.Range("F14").Value = 2317
$ws.Range("F15").Value = 2042
$ws.Range("F16").Value = 717
$ws.Range("F17").Value = 6481
$ws.Range("F18").Value = 123
$ws.Range("F19").Value = 530
$ws.Range("F20").Value = 1228
$ws.Range("F21").Value = 570
$ws.Range("F22").Value = 1477
$ws.Range("F23").Value = 1331
$ws.Range("F24").Value = 1194
$ws.Range("F25").Value = 2296
$ws.Range("F26").Value = 356
$ws.Range("F27").Value = 81
$ws.Range("F29").Value = 731
$ws.Range("F30").Value = 242
$ws.Range("F31").Value = 5303
$ws.Range("F32").Value = 286
$ws.Range("F33").Value = 1254
$ws.Range("F34").Value = 3711
$ws.Range("F35").Value = 301
$ws.Range("F36").Value = 1687
$ws.Range("F38").Value = 159
$ws.Range("F40").Value = 392
$ws.Range("F41").Value = 1768
$ws.Range("F43").Value = 46
$ws.Range("G43").Value = 149
$ws.Range("F44").Value = 105
$ws.Range("F45").Value = 899
$ws.Range("F47").Value = 513
$ws.Range("F48").Value = 197
$ws.Range("F49").Value = 197
$ws.Range("F50").Value = 60
$ws.Range("F51").Value = 78
